# Updates the "Bill Summary" report sheet with refreshed measurement-book
# figures: short point (instead of medium point) quantities, recomputed
# totals, and a new "Tender Premium" / "NET PAYABLE AMOUNT" footer block.
# The report now ends at row 18 instead of row 20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must be stored as TEXT (t="str"/shared-string),
# even when it looks like a number (e.g. "2", "0.00", or an empty string).
# A leading apostrophe forces Excel to treat the input as text instead of
# coercing it to a number; resetting the Style back to "Normal" afterwards
# drops the transient "quote prefix" cell style so only the value changes.
function Set-TextValue($sheet, $addr, $val) {
    $sheet.Range($addr).Value = "'" + $val
    $sheet.Range($addr).Style = "Normal"
}

# --- Row 8: quantity executed upto date changes ---
$ws.Range("C8").Value = 93

# --- Row 9: "Medium point" becomes "Short point" with new quantities ---
$ws.Range("C9").Value = 87
Set-TextValue $ws "D9" "2"
$ws.Range("E9").Value = "Short point (up to 3 mtr.)"
$ws.Range("F9").Value = 256
Set-TextValue $ws "G9" "22272.00"

# --- Row 10: "Long point" quantity/amount recomputed ---
$ws.Range("C10").Value = 22
Set-TextValue $ws "G10" "14564.00"

# --- Row 11: quantity executed upto date changes ---
$ws.Range("C11").Value = 93

# --- Row 12: was "P. point / On board"; now the blank "Total" row ---
Set-TextValue $ws "A12" ""
$ws.Range("C12").Value = 44
Set-TextValue $ws "D12" "8"
$ws.Range("E12").Value = "Total"
$ws.Range("F12").Value = 0
Set-TextValue $ws "G12" "0.00"

# --- Row 13: was the long "P & F ISI marked..." item; now "% Add Tender Premium" ---
Set-TextValue $ws "A13" "%"
$ws.Range("C13").Value = 64
Set-TextValue $ws "D13" "9"
$ws.Range("E13").Value = "Add Tender Premium "
$ws.Range("F13").Value = 0
Set-TextValue $ws "G13" "0.00"

# --- Row 14: was "Total"; now "Grand Total" ---
$ws.Range("C14").Value = 47
Set-TextValue $ws "D14" "10"
$ws.Range("E14").Value = "Grand Total"

# --- Row 15: was "% Add Tender Premium" row; now collapses to a blank A15 ---
$ws.Range("B15:I15").ClearContents()
Set-TextValue $ws "A15" ""

# --- Row 16: was "Grand Total"; now the "Grand Total Rs." totals row ---
Set-TextValue $ws "B16" ""
Set-TextValue $ws "C16" ""
Set-TextValue $ws "D16" ""
$ws.Range("E16").Value = "Grand Total Rs."
Set-TextValue $ws "F16" ""
Set-TextValue $ws "G16" "36836.00"
Set-TextValue $ws "H16" "36836.00"

# --- Row 17: new "Tender Premium @ 0%" totals row ---
Set-TextValue $ws "B17" ""
Set-TextValue $ws "C17" ""
Set-TextValue $ws "D17" ""
Set-TextValue $ws "E17" "Tender Premium @ 0%"
Set-TextValue $ws "F17" ""
Set-TextValue $ws "G17" "0.00"
Set-TextValue $ws "H17" "0.00"
Set-TextValue $ws "I17" ""

# --- Row 18: was "Grand Total Rs."; now "NET PAYABLE AMOUNT Rs." ---
$ws.Range("E18").Value = "NET PAYABLE AMOUNT Rs."
Set-TextValue $ws "G18" "36836.00"
Set-TextValue $ws "H18" "36836.00"

# --- Old rows 19-20 (Tender Premium / NET PAYABLE AMOUNT) no longer needed ---
$ws.Rows("19:20").Delete()
